$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts 3 new rows of fresh data right above the old
# row 18, pushing all the old rows 18-33 down to rows 21-36 (same content,
# just shifted). Insert 3 rows first so everything below moves down intact,
# then populate the 3 freshly-inserted rows with the new week's data.
$ws.Rows("18:20").Insert()

# --- Row 18: Espárragos, Sin especificar / Banquete ---
$ws.Cells.Item(18, 1).Value  = 12
$ws.Cells.Item(18, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(18, 3).Value  = "Metropolitana"
$ws.Cells.Item(18, 4).Value  = 44468
$ws.Cells.Item(18, 5).Value  = 13
$ws.Cells.Item(18, 6).Value  = 300000000
$ws.Cells.Item(18, 7).Value  = "Espárragos"
$ws.Cells.Item(18, 8).Value  = "Sin especificar"
$ws.Cells.Item(18, 9).Value  = "Banquete"
$ws.Cells.Item(18, 10).Value = 90
$ws.Cells.Item(18, 11).Value = 1600
$ws.Cells.Item(18, 12).Value = 1600
$ws.Cells.Item(18, 13).Value = 1600
$ws.Cells.Item(18, 14).Value = "$/kilo"
$ws.Cells.Item(18, 15).Value = "Región Metropolitana"
$ws.Cells.Item(18, 16).Value = 1600
$ws.Cells.Item(18, 17).Value = 1
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# --- Row 19: Espárragos, Sin especificar / Primera ---
$ws.Cells.Item(19, 1).Value  = 12
$ws.Cells.Item(19, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(19, 3).Value  = "Metropolitana"
$ws.Cells.Item(19, 4).Value  = 44468
$ws.Cells.Item(19, 5).Value  = 13
$ws.Cells.Item(19, 6).Value  = 300000000
$ws.Cells.Item(19, 7).Value  = "Espárragos"
$ws.Cells.Item(19, 8).Value  = "Sin especificar"
$ws.Cells.Item(19, 9).Value  = "Primera"
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1500
$ws.Cells.Item(19, 14).Value = "$/kilo"
$ws.Cells.Item(19, 15).Value = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value = 1500
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# --- Row 20: Espárragos, Sin especificar / Segunda ---
$ws.Cells.Item(20, 1).Value  = 12
$ws.Cells.Item(20, 2).Value  = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(20, 3).Value  = "Metropolitana"
$ws.Cells.Item(20, 4).Value  = 44468
$ws.Cells.Item(20, 5).Value  = 13
$ws.Cells.Item(20, 6).Value  = 300000000
$ws.Cells.Item(20, 7).Value  = "Espárragos"
$ws.Cells.Item(20, 8).Value  = "Sin especificar"
$ws.Cells.Item(20, 9).Value  = "Segunda"
$ws.Cells.Item(20, 10).Value = 75
$ws.Cells.Item(20, 11).Value = 1300
$ws.Cells.Item(20, 12).Value = 1300
$ws.Cells.Item(20, 13).Value = 1300
$ws.Cells.Item(20, 14).Value = "$/kilo"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1300
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
